$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparison")

# Delete the three "wardrobe" comparison rows (bottom-up so row numbers stay valid)
$ws.Rows.Item(57).Delete()
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(25).Delete()

# Trim stray leading/trailing spaces from the concept columns and fix three
# mistranslated / unpluralized words.
for ($r = 2; $r -le 66; $r++) {
    $a = $ws.Cells.Item($r, 1).Text
    $b = $ws.Cells.Item($r, 2).Text

    $a = $a.Trim()
    $b = $b.Trim()

    if ($b -eq "pijama") { $b = "pajamas" }
    if ($b -eq "sandal") { $b = "sandals" }
    if ($b -eq "slipper") { $b = "slippers" }

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}
